$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column D holds values that look numeric (e.g. "29.161.44", "1.007")
# but must stay TEXT, exactly as in the source file (t="inlineStr").
# Prefixing with an apostrophe forces Excel to store the value as text
# instead of auto-converting it to a number; ClearFormats() afterwards
# drops the transient "quote prefix" cell style so no stray style index
# is introduced (matches the diff, which shows no "s=" attribute change).
function Set-TextValue($cell, $val) {
    $cell.Value = "'" + $val
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "29.161.44"
$ws.Range("E2").Value = "  +0.20%  "
Set-TextValue $ws.Range("D3") "1.836.44"
$ws.Range("E3").Value = "  +0.14%  "
Set-TextValue $ws.Range("D4") "1.007"
$ws.Range("E4").Value = "  +0.52%  "
Set-TextValue $ws.Range("D5") "243.31"
$ws.Range("E5").Value = "  -0.34%  "
Set-TextValue $ws.Range("D6") "0.6171"
$ws.Range("E6").Value = "  -1.74%  "
Set-TextValue $ws.Range("D7") "1.009"
$ws.Range("E7").Value = "  +0.57%  "
Set-TextValue $ws.Range("D8") "0.07426"
$ws.Range("E8").Value = "  -0.65%  "
Set-TextValue $ws.Range("D9") "0.2919"
$ws.Range("E9").Value = "  -0.20%  "
Set-TextValue $ws.Range("D10") "22.89"
$ws.Range("E10").Value = "  -0.78%  "
Set-TextValue $ws.Range("D11") "0.07705"
$ws.Range("E11").Value = "  -0.30%  "
Set-TextValue $ws.Range("D12") "1.825.79"
$ws.Range("E12").Value = "  -0.54%  "
Set-TextValue $ws.Range("D13") "4.968"
$ws.Range("E13").Value = "  -0.13%  "
Set-TextValue $ws.Range("D14") "0.6686"
$ws.Range("E14").Value = "  +0.05%  "
Set-TextValue $ws.Range("D15") "82.59"
$ws.Range("E15").Value = "  -0.17%  "
Set-TextValue $ws.Range("D16") "0.000009074"
$ws.Range("E16").Value = "  -2.50%  "
Set-TextValue $ws.Range("D17") "5.864"
$ws.Range("E17").Value = "  -2.81%  "
Set-TextValue $ws.Range("D18") "29.160.72"
$ws.Range("E18").Value = "  +0.10%  "
Set-TextValue $ws.Range("D19") "2.084.88"
Set-TextValue $ws.Range("D20") "235.68"
$ws.Range("E20").Value = "  +5.55%  "
Set-TextValue $ws.Range("D21") "12.55"
$ws.Range("E21").Value = "  -0.34%  "
Set-TextValue $ws.Range("D22") "1.009"
$ws.Range("E22").Value = "  +0.44%  "
Set-TextValue $ws.Range("D23") "7.141"
$ws.Range("E23").Value = "  +0.17%  "
Set-TextValue $ws.Range("D24") "1.009"
$ws.Range("E24").Value = "  +0.59%  "
Set-TextValue $ws.Range("D25") "159.24"
$ws.Range("E25").Value = "  -0.57%  "
Set-TextValue $ws.Range("D26") "0.1413"
$ws.Range("E26").Value = "  +1.07%  "
Set-TextValue $ws.Range("D27") "8.483"
$ws.Range("E27").Value = "  -0.22%  "
Set-TextValue $ws.Range("D28") "17.77"
$ws.Range("E28").Value = "  -0.82%  "
Set-TextValue $ws.Range("D29") "1.501"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D30") "4.113"
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D31") "0.05553"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D32") "4.119"
$ws.Range("E32").Value = "  -0.65%  "
Set-TextValue $ws.Range("D33") "1.215"
$ws.Range("E33").Value = "  +0.63%  "
Set-TextValue $ws.Range("D34") "1.839"
$ws.Range("E34").Value = "  -0.70%  "
Set-TextValue $ws.Range("D35") "0.7407"
$ws.Range("E35").Value = "  -1.31%  "
Set-TextValue $ws.Range("D36") "1.136"
$ws.Range("E36").Value = "  -0.08%  "
Set-TextValue $ws.Range("D37") "2.650"
$ws.Range("E37").Value = "  +1.47%  "
Set-TextValue $ws.Range("D38") "2.822"
$ws.Range("E38").Value = "  +2.65%  "
Set-TextValue $ws.Range("D39") "0.01774"
$ws.Range("E39").Value = "  -0.53%  "
Set-TextValue $ws.Range("D40") "1.203.66"
$ws.Range("E40").Value = "  -2.29%  "
Set-TextValue $ws.Range("D41") "6.411"
$ws.Range("E41").Value = "  -2.38%  "
Set-TextValue $ws.Range("D42") "0.8998"
$ws.Range("E42").Value = "  +0.45%  "
Set-TextValue $ws.Range("D43") "1.009"
$ws.Range("E43").Value = "  +0.58%  "
Set-TextValue $ws.Range("D44") "101.43"
$ws.Range("E44").Value = "  -0.45%  "
Set-TextValue $ws.Range("D45") "1.981.73"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D46") "64.81"
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D47") "0.00000000121"
$ws.Range("E47").Value = "  -1.14%  "
Set-TextValue $ws.Range("D48") "0.5123"
$ws.Range("E48").Value = "  +0.64%  "
Set-TextValue $ws.Range("D49") "0.4024"
Set-TextValue $ws.Range("D50") "9.114"
$ws.Range("E50").Value = "  +0.84%  "
Set-TextValue $ws.Range("D51") "0.05822"
$ws.Range("E51").Value = "  +0.25%  "
